# Update "want-to-go" attendee counts (column F) across the worksheets, as
# published by the latest data refresh (commit 456a3b4).
#
# Sheets (by position, matching workbook.xml order):
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Performances)
#   3 = 本地生活 (Local life) - unaffected
#   4 = 全部类型 (All types, a combined/union view)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F5").Value  = 15621
$wsExpo.Range("F7").Value  = 9
$wsExpo.Range("F8").Value  = 705
$wsExpo.Range("F9").Value  = 15417
$wsExpo.Range("F11").Value = 9018
$wsExpo.Range("F12").Value = 383
$wsExpo.Range("F13").Value = 8
$wsExpo.Range("F17").Value = 428
$wsExpo.Range("F20").Value = 51
$wsExpo.Range("F25").Value = 1113
$wsExpo.Range("F26").Value = 1
$wsExpo.Range("F27").Value = 15
$wsExpo.Range("F33").Value = 62
$wsExpo.Range("F34").Value = 41
$wsExpo.Range("F36").Value = 323
$wsExpo.Range("F39").Value = 5540

# --- Sheet 2: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Range("F4").Value = 1

# --- Sheet 4: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value  = 15621
$wsAll.Range("F7").Value  = 9
$wsAll.Range("F8").Value  = 705
$wsAll.Range("F9").Value  = 15417
$wsAll.Range("F11").Value = 9018
$wsAll.Range("F12").Value = 383
$wsAll.Range("F13").Value = 8
$wsAll.Range("F17").Value = 428
$wsAll.Range("F20").Value = 51
$wsAll.Range("F25").Value = 1113
$wsAll.Range("F26").Value = 1
$wsAll.Range("F27").Value = 15
$wsAll.Range("F35").Value = 62
$wsAll.Range("F36").Value = 41
$wsAll.Range("F38").Value = 323
$wsAll.Range("F41").Value = 5540
$wsAll.Range("F42").Value = 1
